# Weekly update: insert a new Mango price record at the top of the
# Vega Monumental Concepción data block (row 103), pushing the existing
# rows (old 103-133) down by one (new 104-134).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row before the current row 103 - everything below
# (including formatting) shifts down by one row automatically.
$ws.Rows.Item(103).Insert()

# Fill the newly inserted row 103 with the new weekly record.
$ws.Range("A103").Value = 11
$ws.Range("B103").Value = "Vega Monumental Concepción"
$ws.Range("C103").Value = "Bíobío"
$ws.Range("D103").Value = 44841
$ws.Range("E103").Value = 8
$ws.Range("F103").Value = "Fruta"
$ws.Range("G103").Value = 100108
$ws.Range("H103").Value = "Tropicales y subtropicales"
$ws.Range("I103").Value = 100108002
$ws.Range("J103").Value = "Mango"
$ws.Range("K103").Value = "Sin especificar"
$ws.Range("L103").Value = "Primera"
$ws.Range("M103").Value = 270
$ws.Range("N103").Value = 7000
$ws.Range("O103").Value = 7500
$ws.Range("P103").Value = 7278
$ws.Range("Q103").Value = "$/bandeja 4 kilos"
$ws.Range("R103").Value = "Brasil"
$ws.Range("S103").Value = 1820
$ws.Range("T103").Value = 4
